# B1--and-B2-PowerPoint.pptx edit
# 1) Slide 5's table gets a new (built-in) table style GUID.
# 2) The presentation's theme colour scheme (bound to the slide master,
#    i.e. ppt/theme/theme2.xml) is repainted from the "Integral / Red
#    Violet" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{086492F8-2A8F-4733-9FDD-C57E13C076E4}")
    }
}

# --- 2. Theme colour scheme -> "Office" palette ---------------------------
function Get-ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officePalette = @{
    1  = "000000"  # Dark 1
    2  = "FFFFFF"  # Light 1
    3  = "44546A"  # Dark 2
    4  = "E7E6E6"  # Light 2
    5  = "5B9BD5"  # Accent 1
    6  = "ED7D31"  # Accent 2
    7  = "A5A5A5"  # Accent 3
    8  = "FFC000"  # Accent 4
    9  = "4472C4"  # Accent 5
    10 = "70AD47"  # Accent 6
    11 = "0563C1"  # Hyperlink
    12 = "954F72"  # Followed hyperlink
}

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($idx in $officePalette.Keys) {
    $colorScheme.Colors($idx).RGB = Get-ComRGB $officePalette[$idx]
}
